# Update the benchmark-results table (1 column x 46 rows) so that it
# reflects the corrected README.md stats.
#
# The first dozen "summary" rows get new values and the last three
# "raw run" rows - which previously held a full tab-separated line of
# per-iteration numbers - get collapsed down to the single summary
# value that used to live in the first three rows.

$d = $word.ActiveDocument
$t = $d.Tables(1)

# Simple single-value replacements (row -> new text)
$t.Cell(1, 1).Range.Text  = "0M"
$t.Cell(2, 1).Range.Text  = "0M"
$t.Cell(3, 1).Range.Text  = "0M"
$t.Cell(4, 1).Range.Text  = "3367"
$t.Cell(5, 1).Range.Text  = "0.00001"
$t.Cell(6, 1).Range.Text  = "0.00255"
$t.Cell(7, 1).Range.Text  = "0.00012"
$t.Cell(8, 1).Range.Text  = "0.00004"
$t.Cell(9, 1).Range.Text  = "0.00018"
$t.Cell(10, 1).Range.Text = "0.00020"
$t.Cell(11, 1).Range.Text = "0.00025"
$t.Cell(12, 1).Range.Text = "0.44210"

# Rows that previously contained a full tab-separated run of numbers
# are collapsed down to a single value.
$t.Cell(44, 1).Range.Text = "99.87"
$t.Cell(45, 1).Range.Text = "0.44"
$t.Cell(46, 1).Range.Text = "337"
